$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new inventory row (row 30) mirroring the existing rows' layout.
$ws.Range("A30").Value = "ZL4VQJ"
$ws.Range("B30").Value = "Chip Epson"
$ws.Range("C30").Value = "F170"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 100000
$ws.Range("F30").Value = 18
$ws.Range("G30").Value = 0
$ws.Range("H30").Formula = "=(E30-D30)*G30"
$ws.Range("I30").Formula = "=D30*F30"
$ws.Range("J30").Value = 0
